$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 5.004338333333333
$ws.Range("H2").Value = 15.013015
$ws.Range("I2").Value = 0.09842252509832523
$ws.Range("J2").Value = 0.09842252509832525
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.05256533333333333
$ws.Range("N2").Value = 0.157696
$ws.Range("O2").Value = 0.7714955259952154
$ws.Range("P2").Value = 0.7714955259952153
$ws.Range("Q2").Value = 0.2630547126044444
$ws.Range("R2").Value = 2.36749241344
$ws.Range("S2").Value = 0.07593253777050972
$ws.Range("T2").Value = 0.07593253777050972

# Row 3
$ws.Range("G3").Value = 5.004338333333333
$ws.Range("H3").Value = 15.013015
$ws.Range("I3").Value = 0.09842252509832523
$ws.Range("J3").Value = 0.09842252509832525
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.3333333333333333
$ws.Range("M3").Value = 0.015569
$ws.Range("N3").Value = 0.046707
$ws.Range("O3").Value = 0.2285044740047847
$ws.Range("P3").Value = 0.2285044740047847
$ws.Range("Q3").Value = 0.07791254351166665
$ws.Range("R3").Value = 0.701212891605
$ws.Range("S3").Value = 0.02248998732781553
$ws.Range("T3").Value = 0.02248998732781553

# Row 4
$ws.Range("I4").Value = 0.02349487489124971
$ws.Range("J4").Value = 0.02349487489124972
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.05256533333333333
$ws.Range("N4").Value = 0.157696
$ws.Range("O4").Value = 0.7714955259952154
$ws.Range("P4").Value = 0.7714955259952153
$ws.Range("Q4").Value = 0.06279495020088889
$ws.Range("R4").Value = 0.565154551808
$ws.Range("S4").Value = 0.01812619086241647
$ws.Range("T4").Value = 0.01812619086241648

# Row 5
$ws.Range("I5").Value = 0.02349487489124971
$ws.Range("J5").Value = 0.02349487489124972
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.015569
$ws.Range("N5").Value = 0.046707
$ws.Range("O5").Value = 0.2285044740047847
$ws.Range("P5").Value = 0.2285044740047847
$ws.Range("Q5").Value = 0.01859884676233333
$ws.Range("R5").Value = 0.167389620861
$ws.Range("S5").Value = 0.005368684028833237
$ws.Range("T5").Value = 0.005368684028833238

# Row 6
$ws.Range("G6").Value = 7.506370333333334
$ws.Range("H6").Value = 22.519111
$ws.Range("I6").Value = 0.1476310899302686
$ws.Range("J6").Value = 0.1476310899302687
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.3333333333333333
$ws.Range("M6").Value = 0.05256533333333333
$ws.Range("N6").Value = 0.157696
$ws.Range("O6").Value = 0.7714955259952154
$ws.Range("P6").Value = 0.7714955259952153
$ws.Range("Q6").Value = 0.3945748586951111
$ws.Range("R6").Value = 3.551173728256
$ws.Range("S6").Value = 0.1138967253789995
$ws.Range("T6").Value = 0.1138967253789996

# Row 7
$ws.Range("G7").Value = 7.506370333333334
$ws.Range("H7").Value = 22.519111
$ws.Range("I7").Value = 0.1476310899302686
$ws.Range("J7").Value = 0.1476310899302687
$ws.Range("K7").Value = 1
$ws.Range("L7").Value = 0.3333333333333333
$ws.Range("M7").Value = 0.015569
$ws.Range("N7").Value = 0.046707
$ws.Range("O7").Value = 0.2285044740047847
$ws.Range("P7").Value = 0.2285044740047847
$ws.Range("Q7").Value = 0.1168666797196667
$ws.Range("R7").Value = 1.051800117477
$ws.Range("S7").Value = 0.0337343645512691
$ws.Range("T7").Value = 0.03373436455126911

# Row 8
$ws.Range("G8").Value = 0.209799
$ws.Range("H8").Value = 0.6293970000000001
$ws.Range("I8").Value = 0.004126209294356305
$ws.Range("J8").Value = 0.004126209294356305
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.05256533333333333
$ws.Range("N8").Value = 0.157696
$ws.Range("O8").Value = 0.7714955259952154
$ws.Range("P8").Value = 0.7714955259952153
$ws.Range("Q8").Value = 0.011028154368
$ws.Range("R8").Value = 0.09925338931200002
$ws.Range("S8").Value = 0.003183352009915764
$ws.Range("T8").Value = 0.003183352009915764

# Row 9
$ws.Range("G9").Value = 0.209799
$ws.Range("H9").Value = 0.6293970000000001
$ws.Range("I9").Value = 0.004126209294356305
$ws.Range("J9").Value = 0.004126209294356305
$ws.Range("K9").Value = 1
$ws.Range("L9").Value = 0.3333333333333333
$ws.Range("M9").Value = 0.015569
$ws.Range("N9").Value = 0.046707
$ws.Range("O9").Value = 0.2285044740047847
$ws.Range("P9").Value = 0.2285044740047847
$ws.Range("Q9").Value = 0.003266360631
$ws.Range("R9").Value = 0.029397245679
$ws.Range("S9").Value = 0.0009428572844405412
$ws.Range("T9").Value = 0.0009428572844405413

# Row 10
$ws.Range("G10").Value = 31.40248233333334
$ws.Range("H10").Value = 94.207447
$ws.Range("I10").Value = 0.6176064445953491
$ws.Range("J10").Value = 0.6176064445953492
$ws.Range("K10").Value = 1
$ws.Range("L10").Value = 0.3333333333333333
$ws.Range("M10").Value = 0.05256533333333333
$ws.Range("N10").Value = 0.157696
$ws.Range("O10").Value = 0.7714955259952154
$ws.Range("P10").Value = 0.7714955259952153
$ws.Range("Q10").Value = 1.650681951345778
$ws.Range("R10").Value = 14.856137562112
$ws.Range("S10").Value = 0.4764806088311237
$ws.Range("T10").Value = 0.4764806088311238

# Row 11
$ws.Range("G11").Value = 31.40248233333334
$ws.Range("H11").Value = 94.207447
$ws.Range("I11").Value = 0.6176064445953491
$ws.Range("J11").Value = 0.6176064445953492
$ws.Range("K11").Value = 1
$ws.Range("L11").Value = 0.3333333333333333
$ws.Range("M11").Value = 0.015569
$ws.Range("N11").Value = 0.046707
$ws.Range("O11").Value = 0.2285044740047847
$ws.Range("P11").Value = 0.2285044740047847
$ws.Range("Q11").Value = 0.4889052474476667
$ws.Range("R11").Value = 4.400147227029
$ws.Range("S11").Value = 0.1411258357642254
$ws.Range("T11").Value = 0.1411258357642254

# Row 12
$ws.Range("G12").Value = 5.52786
$ws.Range("H12").Value = 16.58358
$ws.Range("I12").Value = 0.108718856190451
$ws.Range("J12").Value = 0.1087188561904511
$ws.Range("K12").Value = 1
$ws.Range("L12").Value = 0.3333333333333333
$ws.Range("M12").Value = 0.05256533333333333
$ws.Range("N12").Value = 0.157696
$ws.Range("O12").Value = 0.7714955259952154
$ws.Range("P12").Value = 0.7714955259952153
$ws.Range("Q12").Value = 0.29057380352
$ws.Range("R12").Value = 2.61516423168
$ws.Range("S12").Value = 0.0838761111422502
$ws.Range("T12").Value = 0.08387611114225021

# Row 13
$ws.Range("G13").Value = 5.52786
$ws.Range("H13").Value = 16.58358
$ws.Range("I13").Value = 0.108718856190451
$ws.Range("J13").Value = 0.1087188561904511
$ws.Range("K13").Value = 1
$ws.Range("L13").Value = 0.3333333333333333
$ws.Range("M13").Value = 0.015569
$ws.Range("N13").Value = 0.046707
$ws.Range("O13").Value = 0.2285044740047847
$ws.Range("P13").Value = 0.2285044740047847
$ws.Range("Q13").Value = 0.08606325233999999
$ws.Range("R13").Value = 0.7745692710599998
$ws.Range("S13").Value = 0.02484274504820084
$ws.Range("T13").Value = 0.02484274504820085
